$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 446, shifting the existing
# rows 446..475 down to 448..477 (dimension grows from A1:T475 to A1:T477).
$ws.Rows.Item(446).Insert()
$ws.Rows.Item(446).Insert()

# --- New row 446 ---
$ws.Cells.Item(446, 1).Value = 3
$ws.Cells.Item(446, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(446, 3).Value = "Coquimbo"
$ws.Cells.Item(446, 4).Value = 44826
$ws.Cells.Item(446, 5).Value = 5
$ws.Cells.Item(446, 6).Value = "Fruta"
$ws.Cells.Item(446, 7).Value = 100108
$ws.Cells.Item(446, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(446, 9).Value = 100108002
$ws.Cells.Item(446, 10).Value = "Mango"
$ws.Cells.Item(446, 11).Value = "Sin especificar"
$ws.Cells.Item(446, 12).Value = "Primera"
$ws.Cells.Item(446, 13).Value = 228
$ws.Cells.Item(446, 14).Value = 11000
$ws.Cells.Item(446, 15).Value = 11000
$ws.Cells.Item(446, 16).Value = 11000
$ws.Cells.Item(446, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(446, 18).Value = "Brasil"
$ws.Cells.Item(446, 19).Value = 2750
$ws.Cells.Item(446, 20).Value = 4

# --- New row 447 ---
$ws.Cells.Item(447, 1).Value = 3
$ws.Cells.Item(447, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(447, 3).Value = "Coquimbo"
$ws.Cells.Item(447, 4).Value = 44826
$ws.Cells.Item(447, 5).Value = 5
$ws.Cells.Item(447, 6).Value = "Fruta"
$ws.Cells.Item(447, 7).Value = 100108
$ws.Cells.Item(447, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(447, 9).Value = 100108002
$ws.Cells.Item(447, 10).Value = "Mango"
$ws.Cells.Item(447, 11).Value = "Sin especificar"
$ws.Cells.Item(447, 12).Value = "Segunda"
$ws.Cells.Item(447, 13).Value = 228
$ws.Cells.Item(447, 14).Value = 11000
$ws.Cells.Item(447, 15).Value = 11000
$ws.Cells.Item(447, 16).Value = 11000
$ws.Cells.Item(447, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(447, 18).Value = "Brasil"
$ws.Cells.Item(447, 19).Value = 2750
$ws.Cells.Item(447, 20).Value = 4
